$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2722696666666667
$ws.Range("H2").Value = 0.816809
$ws.Range("I2").Value = 0.01577089893809228
$ws.Range("J2").Value = 0.01577089893809228
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2615913333333333
$ws.Range("N2").Value = 0.784774
$ws.Range("O2").Value = 0.08239613548481725
$ws.Range("P2").Value = 0.08239613548481727
$ws.Range("Q2").Value = 0.07122338512955556
$ws.Range("R2").Value = 0.641010466166
$ws.Range("S2").Value = 0.001299461125620412
$ws.Range("T2").Value = 0.001299461125620412

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2722696666666667
$ws.Range("H3").Value = 0.816809
$ws.Range("I3").Value = 0.01577089893809228
$ws.Range("J3").Value = 0.01577089893809228
$ws.Range("N3").Value = 5.233242000000001
$ws.Range("O3").Value = 0.5494561706387266
$ws.Range("P3").Value = 0.5494561706387268
$ws.Range("Q3").Value = 0.4749510183086668
$ws.Range("R3").Value = 4.274559164778
$ws.Range("S3").Value = 0.008665417738054544
$ws.Range("T3").Value = 0.008665417738054544

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2722696666666667
$ws.Range("H4").Value = 0.816809
$ws.Range("I4").Value = 0.01577089893809228
$ws.Range("J4").Value = 0.01577089893809228
$ws.Range("M4").Value = 1.168795666666667
$ws.Range("N4").Value = 3.506387
$ws.Range("O4").Value = 0.3681476938764561
$ws.Range("P4").Value = 0.3681476938764561
$ws.Range("Q4").Value = 0.3182276065647778
$ws.Range("R4").Value = 2.864048459083
$ws.Range("S4").Value = 0.005806020074417322
$ws.Range("T4").Value = 0.005806020074417322

# Row 5
$ws.Range("I5").Value = 0.8050543166133334
$ws.Range("J5").Value = 0.8050543166133333
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2615913333333333
$ws.Range("N5").Value = 0.784774
$ws.Range("O5").Value = 0.08239613548481725
$ws.Range("P5").Value = 0.08239613548481727
$ws.Range("Q5").Value = 3.635727669516
$ws.Range("R5").Value = 32.721549025644
$ws.Range("S5").Value = 0.06633336454430919
$ws.Range("T5").Value = 0.06633336454430919

# Row 6
$ws.Range("I6").Value = 0.8050543166133334
$ws.Range("J6").Value = 0.8050543166133333
$ws.Range("N6").Value = 5.233242000000001
$ws.Range("O6").Value = 0.5494561706387266
$ws.Range("P6").Value = 0.5494561706387268
$ws.Range("S6").Value = 0.4423420619625392
$ws.Range("T6").Value = 0.4423420619625392

# Row 7
$ws.Range("I7").Value = 0.8050543166133334
$ws.Range("J7").Value = 0.8050543166133333
$ws.Range("M7").Value = 1.168795666666667
$ws.Range("N7").Value = 3.506387
$ws.Range("O7").Value = 0.3681476938764561
$ws.Range("P7").Value = 0.3681476938764561
$ws.Range("Q7").Value = 16.244508910758
$ws.Range("R7").Value = 146.200580196822
$ws.Range("S7").Value = 0.296378890106485
$ws.Range("T7").Value = 0.296378890106485

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.093283333333334
$ws.Range("H8").Value = 9.279850000000001
$ws.Range("I8").Value = 0.1791747844485745
$ws.Range("J8").Value = 0.1791747844485744
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2615913333333333
$ws.Range("N8").Value = 0.784774
$ws.Range("O8").Value = 0.08239613548481725
$ws.Range("P8").Value = 0.08239613548481727
$ws.Range("Q8").Value = 0.8091761115444446
$ws.Range("R8").Value = 7.282585003900001
$ws.Range("S8").Value = 0.01476330981488767
$ws.Range("T8").Value = 0.01476330981488767

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.093283333333334
$ws.Range("H9").Value = 9.279850000000001
$ws.Range("I9").Value = 0.1791747844485745
$ws.Range("J9").Value = 0.1791747844485744
$ws.Range("N9").Value = 5.233242000000001
$ws.Range("O9").Value = 0.5494561706387266
$ws.Range("P9").Value = 0.5494561706387268
$ws.Range("Q9").Value = 5.395966752633335
$ws.Range("R9").Value = 48.56370077370001
$ws.Range("S9").Value = 0.09844869093813299
$ws.Range("T9").Value = 0.09844869093813299

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.093283333333334
$ws.Range("H10").Value = 9.279850000000001
$ws.Range("I10").Value = 0.1791747844485745
$ws.Range("J10").Value = 0.1791747844485744
$ws.Range("M10").Value = 1.168795666666667
$ws.Range("N10").Value = 3.506387
$ws.Range("O10").Value = 0.3681476938764561
$ws.Range("P10").Value = 0.3681476938764561
$ws.Range("Q10").Value = 3.615416155772223
$ws.Range("R10").Value = 32.53874540195001
$ws.Range("S10").Value = 0.06596278369555379
$ws.Range("T10").Value = 0.06596278369555379
